# Fix bugs for stats2: the data rows (2-21) were reshuffled into the
# correct order. Column A = symbol id, columns B:F = reel1..reel5 counts.
# Rewrite the block A2:F21 with the corrected row order (values only,
# same data set, just reordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1001, 18, 30, 75, 60, 72),
    @(201,  9,  30, 15, 45, 30),
    @(1202, 2,  10, 10, 10, 10),
    @(101,  9,  30, 15, 60, 15),
    @(901,  16, 15, 45, 60, 60),
    @(902,  1,  0,  0,  0,  0),
    @(301,  6,  45, 30, 60, 45),
    @(501,  9,  52, 30, 75, 45),
    @(801,  3,  67, 65, 52, 45),
    @(1203, 3,  15, 15, 15, 15),
    @(401,  9,  48, 67, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(601,  9,  60, 67, 60, 42),
    @(1201, 2,  10, 10, 10, 10),
    @(802,  0,  4,  5,  4,  0),
    @(1,    0,  2,  2,  2,  2),
    @(2,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(1101, 0,  15, 30, 30, 0),
    @(3,    0,  3,  3,  3,  3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
